# Update the multiplication-fact table. Cells are addressed by
# (row, column) rather than by searching for text, because some values
# are reused across cells (e.g. "413×8=3304" appears both as a new
# value in row 1 and as an old value in row 5) -- a text-based
# Find/Replace could clobber the wrong occurrence.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Old = "339×3=1017"; New = "815×4=3260" },
    @{ Row = 1;  Col = 2; Old = "546×5=2730"; New = "174×9=1566" },
    @{ Row = 1;  Col = 3; Old = "910×3=2730"; New = "413×8=3304" },
    @{ Row = 1;  Col = 4; Old = "584×6=3504"; New = "625×9=5625" },
    @{ Row = 1;  Col = 5; Old = "750×9=6750"; New = "706×9=6354" },

    @{ Row = 5;  Col = 1; Old = "413×8=3304"; New = "571×5=2855" },
    @{ Row = 5;  Col = 2; Old = "933×9=8397"; New = "521×2=1042" },
    @{ Row = 5;  Col = 3; Old = "789×4=3156"; New = "739×2=1478" },
    @{ Row = 5;  Col = 4; Old = "319×2=638";  New = "741×9=6669" },
    @{ Row = 5;  Col = 5; Old = "518×8=4144"; New = "880×3=2640" },

    @{ Row = 10; Col = 1; Old = "925×7=6475"; New = "285×3=855"  },
    @{ Row = 10; Col = 2; Old = "787×6=4722"; New = "835×8=6680" },
    @{ Row = 10; Col = 3; Old = "424×8=3392"; New = "211×8=1688" },
    @{ Row = 10; Col = 4; Old = "422×4=1688"; New = "327×5=1635" },
    @{ Row = 10; Col = 5; Old = "881×8=7048"; New = "597×7=4179" },

    @{ Row = 15; Col = 1; Old = "805×6=4830"; New = "599×7=4193" },
    @{ Row = 15; Col = 2; Old = "937×9=8433"; New = "159×5=795"  },
    @{ Row = 15; Col = 3; Old = "516×4=2064"; New = "531×6=3186" },
    @{ Row = 15; Col = 4; Old = "233×8=1864"; New = "413×7=2891" },
    @{ Row = 15; Col = 5; Old = "674×7=4718"; New = "743×5=3715" },

    @{ Row = 20; Col = 1; Old = "245×6=1470"; New = "963×4=3852" },
    @{ Row = 20; Col = 2; Old = "366×2=732";  New = "673×7=4711" },
    @{ Row = 20; Col = 3; Old = "505×8=4040"; New = "904×2=1808" },
    @{ Row = 20; Col = 4; Old = "312×4=1248"; New = "133×8=1064" },
    @{ Row = 20; Col = 5; Old = "970×7=6790"; New = "735×7=5145" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cellRange = $cell.Range
    # Cell.Range.Text includes the trailing end-of-cell mark (CR+BEL);
    # strip it before comparing against the expected plain value.
    $current = $cellRange.Text.TrimEnd([char]13, [char]7)
    if ($current -ne $u.Old) {
        Write-Host "WARNING: cell ($($u.Row),$($u.Col)) expected '$($u.Old)' but found '$current'"
    }
    $cellRange.Text = $u.New
}

Write-Host "Done updating $($updates.Count) cells."
